$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.651777982711792
$ws.Range("B1").Value = 1.826148271560669
$ws.Range("C1").Value = 1.888054251670837
$ws.Range("D1").Value = 2.448972702026367
$ws.Range("E1").Value = 3.520786046981812
